$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.003.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.45%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.483.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.78%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'585.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.07%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'168.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -5.57%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.20%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.517"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.13%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.482.03"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.62%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -5.21%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.02%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.338"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.54%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -4.46%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'25.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -4.54%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.936.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.55%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -4.09%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'66.888.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.32%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.478.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.03%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'11.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.82%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.10%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'360.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.13%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.53%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.33%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'4.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -7.08%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'70.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.66%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -6.75%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -8.45%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.22%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.606.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.01%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0₃0935"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -6.84%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'8.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.38%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'507.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -7.57%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.73%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -6.29%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.05%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.127"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.57%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'158.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.28%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -4.13%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.00%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'18.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.81%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -5.15%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -5.95%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.335"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -6.73%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.08%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -4.03%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'39.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.10%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'142.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.45%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.539"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -5.08%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'3.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -4.37%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0₆0268"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -4.98%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.27%  "
$ws.Range("E51").Style = "Normal"
